# Rotate the contents of rows 4-10 up by one, wrapping the first row's
# original content around to the last row:
#   new row 4  = old row 5
#   new row 5  = old row 6
#   new row 6  = old row 7
#   new row 7  = old row 8
#   new row 8  = old row 9
#   new row 9  = old row 10
#   new row 10 = old row 4
#
# Columns Y:AB (Startdatum/Starttid/Slutdatum/Sluttid) hold the exact same
# text in every one of these rows both before and after the edit, so they
# are left untouched here. Excel's COM automation auto-coerces a
# date-formatted string like "2023-06-13" into a real Date value whenever
# it passes back through Range.Value, which would otherwise corrupt those
# text cells for no net change. Splitting the copy into A:X and AC:AY
# sidesteps that round-trip entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 4
$lastRow = 10

# Snapshot every row's two value blocks before mutating anything.
$snapLeft = @{}
$snapRight = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapLeft[$r] = $ws.Range("A" + $r + ":X" + $r).Value()
    $snapRight[$r] = $ws.Range("AC" + $r + ":AY" + $r).Value()
}

# Write row r <- snapshot of row r+1 (for all but the last row).
for ($r = $firstRow; $r -lt $lastRow; $r++) {
    $ws.Range("A" + $r + ":X" + $r).Value = $snapLeft[$r + 1]
    $ws.Range("AC" + $r + ":AY" + $r).Value = $snapRight[$r + 1]
}

# Wrap the very first row's original content around to the last row.
$ws.Range("A" + $lastRow + ":X" + $lastRow).Value = $snapLeft[$firstRow]
$ws.Range("AC" + $lastRow + ":AY" + $lastRow).Value = $snapRight[$firstRow]
